$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Row 6: add "Gestion du projet, documentation et présentation" to E6
$ws.Range("E6").Value = "Gestion du projet, documentation et présentation"

# Row 7: add same text to E7
$ws.Range("E7").Value = "Gestion du projet, documentation et présentation"

# Row 8: fill in Date / Début / Fin / Quoi
$ws.Range("A8").Value = 45475
$ws.Range("B8").Value = 0.76388888888888884
$ws.Range("C8").Value = 0.80555555555555547
$ws.Range("E8").Value = "Gestion du projet, documentation et présentation"

# Row 9: fill in Date / Début / Fin / Quoi
$ws.Range("A9").Value = 45476
$ws.Range("B9").Value = 0.41666666666666669
$ws.Range("C9").Value = 0.52083333333333337
$ws.Range("E9").Value = "Gestion du projet, documentation et présentation"

# Update active selection to G10
$ws.Range("G10").Select()

$wb.Save()
